$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve the existing "header/key-column" cell format (style index 1:
#     bold, centered, top-aligned, thin-bordered) by copying it from A1 before
#     we touch anything, then reapplying it to the new header row + glycan column. ---
$ws.Range("A1").Copy()
$ws.Range("A1:G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2:A5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Clear leftover values/formats in the old used range that fall outside the new A1:G5 footprint ---
$ws.Range("H1:L2").Clear()

# --- Header row (row 1): glycan, binding_score, monosaccharides, motifs, sasa, flexibility, has_multi_node_motifs ---
$headers = @("glycan", "binding_score", "monosaccharides", "motifs", "sasa", "flexibility", "has_multi_node_motifs")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Data rows (rows 2-5) ---
$data = @(
    @("Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)Glc", -0.4710597589339519, "['Gal(b1-4)', 'GlcNAc(b1-3)']", "['Gal(b1-4)GlcNAc(b1-3)']", 5.338285572580087, 0.91117855161729, $true),
    @("Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc", 2.070241755787822, "['Gal(b1-4)', 'GlcNAc(b1-3)']", "['Gal(b1-4)GlcNAc(b1-3)']", 5.27278254643194, 2.180924532303609, $true),
    @("Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc", 2.069056814377114, "['Gal(b1-4)', 'GlcNAc(b1-3)']", "['Gal(b1-4)GlcNAc(b1-3)']", 5.252147263686476, 1.839554809126105, $true),
    @("GlcNAc(b1-4)GlcNAc(b1-4)GlcNAc", 1.382957652431078, "['GlcNAc(b1-4)', 'GlcNAc(b1-4)']", "['GlcNAc(b1-4)GlcNAc']", 5.702079978569953, 0.6713220512263312, $true)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $rowIndex++
}
